# Generate Report for Handoff
# Adds two new source-file rows (3c63a561-... and 5cc063e8-...) to every
# worksheet (Overview, zh-cn, de-de), pushing the existing
# ".localization-config" row down, and records their handoff status /
# handoff file / handoff datetime information.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276   # matches existing custom "HyperLink" style color FF6495ED

function Set-Text($ws, $addr, $text) {
    $ws.Range($addr).Value = $text
}

function Add-Link($ws, $addr, $url, $text) {
    $ws.Range($addr).Value = $text
    $ws.Hyperlinks.Add($ws.Range($addr), $url, $null, $null, $text) | Out-Null
    $ws.Range($addr).Font.Underline = 2
    $ws.Range($addr).Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Remove every existing hyperlink (and their relationships) so that they can
# be re-created in the correct, final order (this engine only supports
# removing the whole collection at once).
$ws.Hyperlinks.Delete()

Add-Link $ws "A2" "https://github.com/OpenLocalizationTest/oltest/blob/f8010f942051cba594ec0f688499eaf9aad5cc0f/e2e/0ff38e90-feda-4738-9db6-d2b46e203e42.md" "0ff38e90-feda-4738-9db6-d2b46e203e42.md"
Set-Text $ws "B2" "In Translation"
Set-Text $ws "C2" "In Translation"

Add-Link $ws "A3" "https://github.com/OpenLocalizationTest/oltest/blob/f8010f942051cba594ec0f688499eaf9aad5cc0f/e2e/5fa07b75-d073-426d-abeb-71d17b3c3c40.md" "5fa07b75-d073-426d-abeb-71d17b3c3c40.md"
Set-Text $ws "B3" "In Translation"
Set-Text $ws "C3" "In Translation"

Add-Link $ws "A4" "https://github.com/OpenLocalizationTest/oltest/blob/f8010f942051cba594ec0f688499eaf9aad5cc0f/e2e/3c63a561-8466-498f-ad47-f592cd24babe.md" "3c63a561-8466-498f-ad47-f592cd24babe.md"
Set-Text $ws "B4" "Ready for handoff"
Set-Text $ws "C4" "Ready for handoff"

Add-Link $ws "A5" "https://github.com/OpenLocalizationTest/oltest/blob/f8010f942051cba594ec0f688499eaf9aad5cc0f/e2e/5cc063e8-c3f6-4418-88f1-4716df275030.md" "5cc063e8-c3f6-4418-88f1-4716df275030.md"
Set-Text $ws "B5" "Ready for handoff"
Set-Text $ws "C5" "Ready for handoff"

Add-Link $ws "A6" "https://github.com/OpenLocalizationTest/oltest/blob/f8010f942051cba594ec0f688499eaf9aad5cc0f/.localization-config" ".localization-config"
Set-Text $ws "B6" "Not to be localized"
Set-Text $ws "C6" "Not to be localized"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()

Add-Link $ws "A2" "https://github.com/OpenLocalizationTest/oltest/blob/f8010f942051cba594ec0f688499eaf9aad5cc0f/e2e/0ff38e90-feda-4738-9db6-d2b46e203e42.md" "0ff38e90-feda-4738-9db6-d2b46e203e42.md"
Set-Text $ws "B2" "In Translation"
Add-Link $ws "C2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dada16466475408bdd8930ab202ec2996af2f76c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/0ff38e90-feda-4738-9db6-d2b46e203e42.4a13de436ee2276670ed4d9fdc9f8fad96dff818.zh-cn.xlf" "0ff38e90-feda-4738-9db6-d2b46e203e42.4a13de436ee2276670ed4d9fdc9f8fad96dff818.zh-cn.xlf"
Set-Text $ws "D2" "2016-03-04 02:56:20"
$ws.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-Text $ws "G2" "0001-01-01 00:00:00"
Set-Text $ws "H2" "Include"

Add-Link $ws "A3" "https://github.com/OpenLocalizationTest/oltest/blob/f8010f942051cba594ec0f688499eaf9aad5cc0f/e2e/5fa07b75-d073-426d-abeb-71d17b3c3c40.md" "5fa07b75-d073-426d-abeb-71d17b3c3c40.md"
Set-Text $ws "B3" "In Translation"
Add-Link $ws "C3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dada16466475408bdd8930ab202ec2996af2f76c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/5fa07b75-d073-426d-abeb-71d17b3c3c40.623056b9397dd87897c07d4bf5f90bcd025ce7cd.zh-cn.xlf" "5fa07b75-d073-426d-abeb-71d17b3c3c40.623056b9397dd87897c07d4bf5f90bcd025ce7cd.zh-cn.xlf"
Set-Text $ws "D3" "2016-03-04 02:56:20"
$ws.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-Text $ws "G3" "0001-01-01 00:00:00"
Set-Text $ws "H3" "Include"

Add-Link $ws "A4" "https://github.com/OpenLocalizationTest/oltest/blob/f8010f942051cba594ec0f688499eaf9aad5cc0f/e2e/3c63a561-8466-498f-ad47-f592cd24babe.md" "3c63a561-8466-498f-ad47-f592cd24babe.md"
Set-Text $ws "B4" "Ready for handoff"
Add-Link $ws "C4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dada16466475408bdd8930ab202ec2996af2f76c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/3c63a561-8466-498f-ad47-f592cd24babe.4f4f7b09c9e792f5318cba88d8dab4c326d53b1a.zh-cn.xlf" "3c63a561-8466-498f-ad47-f592cd24babe.4f4f7b09c9e792f5318cba88d8dab4c326d53b1a.zh-cn.xlf"
Set-Text $ws "D4" "2016-03-04 02:58:12"
$ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-Text $ws "G4" "0001-01-01 00:00:00"
Set-Text $ws "H4" "Include"

Add-Link $ws "A5" "https://github.com/OpenLocalizationTest/oltest/blob/f8010f942051cba594ec0f688499eaf9aad5cc0f/e2e/5cc063e8-c3f6-4418-88f1-4716df275030.md" "5cc063e8-c3f6-4418-88f1-4716df275030.md"
Set-Text $ws "B5" "Ready for handoff"
Add-Link $ws "C5" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dada16466475408bdd8930ab202ec2996af2f76c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/5cc063e8-c3f6-4418-88f1-4716df275030.72682c66d925768ecda4a967f7785ac7212805a9.zh-cn.xlf" "5cc063e8-c3f6-4418-88f1-4716df275030.72682c66d925768ecda4a967f7785ac7212805a9.zh-cn.xlf"
Set-Text $ws "D5" "2016-03-04 02:58:12"
$ws.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-Text $ws "G5" "0001-01-01 00:00:00"
Set-Text $ws "H5" "Include"

Add-Link $ws "A6" "https://github.com/OpenLocalizationTest/oltest/blob/f8010f942051cba594ec0f688499eaf9aad5cc0f/.localization-config" ".localization-config"
Set-Text $ws "B6" "Not to be localized"
Set-Text $ws "D6" "0001-01-01 00:00:00"
$ws.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-Text $ws "G6" "0001-01-01 00:00:00"
Set-Text $ws "H6" "Ignored"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()

Add-Link $ws "A2" "https://github.com/OpenLocalizationTest/oltest/blob/f8010f942051cba594ec0f688499eaf9aad5cc0f/e2e/0ff38e90-feda-4738-9db6-d2b46e203e42.md" "0ff38e90-feda-4738-9db6-d2b46e203e42.md"
Set-Text $ws "B2" "In Translation"
Add-Link $ws "C2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2aff849af8a8ce1ced08b8baedac12f2ad00a39a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/0ff38e90-feda-4738-9db6-d2b46e203e42.4a13de436ee2276670ed4d9fdc9f8fad96dff818.de-de.xlf" "0ff38e90-feda-4738-9db6-d2b46e203e42.4a13de436ee2276670ed4d9fdc9f8fad96dff818.de-de.xlf"
Set-Text $ws "D2" "2016-03-04 02:56:47"
$ws.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-Text $ws "G2" "0001-01-01 00:00:00"
Set-Text $ws "H2" "Include"

Add-Link $ws "A3" "https://github.com/OpenLocalizationTest/oltest/blob/f8010f942051cba594ec0f688499eaf9aad5cc0f/e2e/5fa07b75-d073-426d-abeb-71d17b3c3c40.md" "5fa07b75-d073-426d-abeb-71d17b3c3c40.md"
Set-Text $ws "B3" "In Translation"
Add-Link $ws "C3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2aff849af8a8ce1ced08b8baedac12f2ad00a39a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/5fa07b75-d073-426d-abeb-71d17b3c3c40.623056b9397dd87897c07d4bf5f90bcd025ce7cd.de-de.xlf" "5fa07b75-d073-426d-abeb-71d17b3c3c40.623056b9397dd87897c07d4bf5f90bcd025ce7cd.de-de.xlf"
Set-Text $ws "D3" "2016-03-04 02:56:47"
$ws.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-Text $ws "G3" "0001-01-01 00:00:00"
Set-Text $ws "H3" "Include"

Add-Link $ws "A4" "https://github.com/OpenLocalizationTest/oltest/blob/f8010f942051cba594ec0f688499eaf9aad5cc0f/e2e/3c63a561-8466-498f-ad47-f592cd24babe.md" "3c63a561-8466-498f-ad47-f592cd24babe.md"
Set-Text $ws "B4" "Ready for handoff"
Add-Link $ws "C4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2aff849af8a8ce1ced08b8baedac12f2ad00a39a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/3c63a561-8466-498f-ad47-f592cd24babe.4f4f7b09c9e792f5318cba88d8dab4c326d53b1a.de-de.xlf" "3c63a561-8466-498f-ad47-f592cd24babe.4f4f7b09c9e792f5318cba88d8dab4c326d53b1a.de-de.xlf"
Set-Text $ws "D4" "2016-03-04 02:58:25"
$ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-Text $ws "G4" "0001-01-01 00:00:00"
Set-Text $ws "H4" "Include"

Add-Link $ws "A5" "https://github.com/OpenLocalizationTest/oltest/blob/f8010f942051cba594ec0f688499eaf9aad5cc0f/e2e/5cc063e8-c3f6-4418-88f1-4716df275030.md" "5cc063e8-c3f6-4418-88f1-4716df275030.md"
Set-Text $ws "B5" "Ready for handoff"
Add-Link $ws "C5" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2aff849af8a8ce1ced08b8baedac12f2ad00a39a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/5cc063e8-c3f6-4418-88f1-4716df275030.72682c66d925768ecda4a967f7785ac7212805a9.de-de.xlf" "5cc063e8-c3f6-4418-88f1-4716df275030.72682c66d925768ecda4a967f7785ac7212805a9.de-de.xlf"
Set-Text $ws "D5" "2016-03-04 02:58:25"
$ws.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-Text $ws "G5" "0001-01-01 00:00:00"
Set-Text $ws "H5" "Include"

Add-Link $ws "A6" "https://github.com/OpenLocalizationTest/oltest/blob/f8010f942051cba594ec0f688499eaf9aad5cc0f/.localization-config" ".localization-config"
Set-Text $ws "B6" "Not to be localized"
Set-Text $ws "D6" "0001-01-01 00:00:00"
$ws.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-Text $ws "G6" "0001-01-01 00:00:00"
Set-Text $ws "H6" "Ignored"
